# chore: update Sheets via scheduled runner
# Refresh cached Kraken market-price figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# leve-profit sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 287.6
$ws.Range("I9").Value = 311.33334
$ws.Range("J9").Value = 252
$ws.Range("K9").Value = 311.33334
$ws.Range("L9").Value = 252
$ws.Range("M9").Value = -142.33334
$ws.Range("N9").Value = -590
$ws.Range("H40").Value = 6509.091
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 7060
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 7060
$ws.Range("M40").Value = -825
$ws.Range("N40").Value = -7410
$ws.Range("H47").Value = 9999
$ws.Range("I47").Value = 9999
$ws.Range("K47").Value = 9999
$ws.Range("M47").Value = -9027
$ws.Range("H62").Value = 2252.75
$ws.Range("I62").Value = 3002.5
$ws.Range("K62").Value = 3002.5
$ws.Range("M62").Value = -2378.5
$ws.Range("H65").Value = 2252.75
$ws.Range("I65").Value = 3002.5
$ws.Range("K65").Value = 15012.5
$ws.Range("M65").Value = -11892.5
$ws.Range("H86").Value = 2379
$ws.Range("I86").Value = 916.3333
$ws.Range("J86").Value = 3476
$ws.Range("K86").Value = 916.3333
$ws.Range("L86").Value = 3476
$ws.Range("M86").Value = 206.6667
$ws.Range("N86").Value = -5722
$ws.Range("H89").Value = 2379
$ws.Range("I89").Value = 916.3333
$ws.Range("J89").Value = 3476
$ws.Range("K89").Value = 4581.6665
$ws.Range("L89").Value = 17380
$ws.Range("M89").Value = 1034.3335
$ws.Range("N89").Value = -28612
$ws.Range("H112").Value = 3150
$ws.Range("J112").Value = 3560
$ws.Range("L112").Value = 10680
$ws.Range("N112").Value = -12896
$ws.Range("H133").Value = 99780
$ws.Range("J133").Value = 99780
$ws.Range("L133").Value = 99780
$ws.Range("N133").Value = -109900
$ws.Range("H137").Value = 1407.1875
$ws.Range("I137").Value = 1375.3334
$ws.Range("J137").Value = 1448.1428
$ws.Range("K137").Value = 4126.0002
$ws.Range("L137").Value = 4344.428400000001
$ws.Range("M137").Value = -1576.0002
$ws.Range("N137").Value = -9444.428400000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4196.25
$ws.Range("I2").Value = 4196.25
$ws.Range("K2").Value = 4196.25
$ws.Range("M2").Value = -4083.25
$ws.Range("H116").Value = 4196.25
$ws.Range("I116").Value = 4196.25
$ws.Range("K116").Value = 4196.25
$ws.Range("M116").Value = -1902.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4196.25
$ws.Range("I3").Value = 4196.25
$ws.Range("K3").Value = 4196.25
$ws.Range("M3").Value = -4082.25
$ws.Range("H134").Value = 5905.4614
$ws.Range("I134").Value = 968.7143
$ws.Range("K134").Value = 2906.1429
$ws.Range("M134").Value = -371.1428999999998

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5625281.5
$ws.Range("I6").Value = 2857464.8
$ws.Range("K6").Value = 2857464.8
$ws.Range("M6").Value = -2857351.8
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150
$ws.Range("H86").Value = 4811.75
$ws.Range("I86").Value = 4811.75
$ws.Range("K86").Value = 4811.75
$ws.Range("M86").Value = -3688.75
$ws.Range("H89").Value = 4811.75
$ws.Range("I89").Value = 4811.75
$ws.Range("K89").Value = 24058.75
$ws.Range("M89").Value = -18442.75
$ws.Range("H105").Value = 1344.75
$ws.Range("I105").Value = 1459.6666
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 1459.6666
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 287.3334
$ws.Range("N105").Value = -4494

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 5000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 5000
$ws.Range("K114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("M114").Value = 15000
$ws.Range("N114").Value = -21508
$ws.Range("H117").Value = 1624.75
$ws.Range("I117").Value = 999
$ws.Range("K117").Value = 2997
$ws.Range("M117").Value = 445
$ws.Range("H121").Value = 2239.8
$ws.Range("I121").Value = 999
$ws.Range("J121").Value = 2550
$ws.Range("K121").Value = 2997
$ws.Range("L121").Value = 7650
$ws.Range("M121").Value = -1687
$ws.Range("N121").Value = -10270
$ws.Range("H129").Value = 1130
$ws.Range("I129").Value = 1130
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3390
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = 1610

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4625.25
$ws.Range("I70").Value = 4500
$ws.Range("J70").Value = 5001
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 5001
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -5541
$ws.Range("H73").Value = 4625.25
$ws.Range("I73").Value = 4500
$ws.Range("J73").Value = 5001
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 5001
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -6873
$ws.Range("H80").Value = 5250
$ws.Range("I80").Value = 5250
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 5250
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -4252
$ws.Range("H83").Value = 5250
$ws.Range("I83").Value = 5250
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 26250
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -21258

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3899.4
$ws.Range("I61").Value = 3899.4
$ws.Range("K61").Value = 3899.4
$ws.Range("M61").Value = -3697.4
$ws.Range("H113").Value = 3899.4
$ws.Range("I113").Value = 3899.4
$ws.Range("K113").Value = 3899.4
$ws.Range("M113").Value = -1729.4
$ws.Range("H132").Value = 4437.0835
$ws.Range("I132").Value = 4305.5557
$ws.Range("K132").Value = 12916.6671
$ws.Range("M132").Value = -10386.6671

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3299.5715
$ws.Range("I107").Value = 1819.4
$ws.Range("K107").Value = 5458.200000000001
$ws.Range("M107").Value = -3538.200000000001
